$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.952.10'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '3.568.55'
$ws.Range("E3").Value = '  +2.61%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.57'
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.43'
$ws.Range("E6").Value = '  -1.03%  '
$ws.Range("D7").Value = '3.568.43'
$ws.Range("E7").Value = '  +2.68%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("E10").Value = '  +0.68%  '
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.385'
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").Value = '4.176.78'
$ws.Range("E13").Value = '  +2.72%  '
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("D15").Value = '3.572.40'
$ws.Range("E15").Value = '  +2.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.06'
$ws.Range("E16").Value = '  +2.38%  '
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").Value = '65.101.69'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.00'
$ws.Range("E19").Value = '  +3.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.37'
$ws.Range("E20").Value = '  +3.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.84'
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.86'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("E23").Value = '  +4.65%  '
$ws.Range("D24").Value = '3.713.92'
$ws.Range("E24").Value = '  +2.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.16'
$ws.Range("E25").Value = '  +2.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000115'
$ws.Range("E27").Value = '  +6.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.75'
$ws.Range("E28").Value = '  +5.99%  '
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.29'
$ws.Range("E30").Value = '  +3.49%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.44'
$ws.Range("E31").Value = '  +3.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.48'
$ws.Range("E32").Value = '  +24.70%  '
$ws.Range("D33").Value = '3.571.21'
$ws.Range("E33").Value = '  +2.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.02'
$ws.Range("E34").Value = '  +4.08%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.144'
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.95'
$ws.Range("E37").Value = '  +1.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '169.46'
$ws.Range("E38").Value = '  -1.26%  '
$ws.Range("E39").Value = '  +5.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.99'
$ws.Range("E40").Value = '  +5.51%  '
$ws.Range("E41").Value = '  +4.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.17'
$ws.Range("E42").Value = '  +10.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.825'
$ws.Range("E43").Value = '  +1.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.71'
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.47'
$ws.Range("E46").Value = '  +2.81%  '
$ws.Range("E47").Value = '  +5.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.64'
$ws.Range("E48").Value = '  +1.44%  '
$ws.Range("D49").Value = '2.484.47'
$ws.Range("E49").Value = '  +12.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.90'
$ws.Range("E50").Value = '  +3.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.38'
$ws.Range("E51").Value = '  +10.12%  '
